$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Di4+XRhod-1")
$ws.Activate()

# ---- Column widths for H, I, J ----
$ws.Columns.Item(8).ColumnWidth = 16.140625
$ws.Columns.Item(9).ColumnWidth = 10.5703125
$ws.Columns.Item(10).ColumnWidth = 10.42578125

# ---- Row 28: "Parts" section header ----
$ws.Range("H28").Value = "Parts"
$ws.Range("H28:J28").Borders.LineStyle = 1
$ws.Range("H28").Font.Bold = $true

# ---- Row 29: Dichroic mirror ----
$ws.Range("H29").Value = "Dichroic mirror"
$ws.Range("I29").Value = 250
$ws.Range("J29").Value = 250
$ws.Range("H29:J29").Borders.LineStyle = 1

# ---- Row 30: Emission Filter ----
$ws.Range("H30").Value = "Emission Filter"
$ws.Range("I30").Value = 365
$ws.Range("J30").Value = 365
$ws.Range("H30:J30").Borders.LineStyle = 1

# ---- Row 31: Filter Cube ----
$ws.Range("H31").Value = "Filter Cube"
$ws.Range("I31").Value = 494
$ws.Range("J31").Value = "-"
$ws.Range("H31:J31").Borders.LineStyle = 1
$ws.Range("J31").HorizontalAlignment = -4108
$ws.Range("J31").VerticalAlignment = -4108

# ---- Row 32: "Dyes" section header ----
$ws.Range("H32").Value = "Dyes"
$ws.Range("H32:J32").Borders.LineStyle = 1
$ws.Range("H32").Font.Bold = $true

# ---- Row 33: X-Rhod-1 ----
$ws.Range("H33").Value = "X-Rhod-1"
$ws.Range("I33").Value = 454
$ws.Range("J33").Value = 454
$ws.Range("H33:J33").Borders.LineStyle = 1

# ---- Row 34: Cal-630 ----
$ws.Range("H34").Value = "Cal-630"
$ws.Range("I34").Value = 341
$ws.Range("J34").Value = 341
$ws.Range("H34:J34").Borders.LineStyle = 1

# ---- Row 35: Total ----
$ws.Range("H35").Value = "Total"
$ws.Range("I35").Formula = "=SUM(I29:I34)"
$ws.Range("J35").Formula = "=SUM(J29:J34)"
$ws.Range("H35:J35").Borders.LineStyle = 1
$ws.Range("H35").Font.Bold = $true
$ws.Range("I35:J35").Font.Bold = $true

# ---- Currency format for the numeric columns ----
$currencyFmt = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'
$ws.Range("I29:J30").NumberFormat = $currencyFmt
$ws.Range("I31").NumberFormat = $currencyFmt
$ws.Range("J31").NumberFormat = $currencyFmt
$ws.Range("I33:J35").NumberFormat = $currencyFmt

# ---- Page setup ----
$ws.PageSetup.Orientation = 1

# ---- Selection ----
$ws.Range("L31").Select()

Write-Output "done"
